$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New query text for the "Remove the door" (row 6) and
# "Delete all walls that creates with Base Wall a room" (row 7) entries.
# Fill order matters: it drives the order new strings are appended to the
# shared-strings table, so E7 is set first, then D6, E6, D7.
$e7 = @'
MATCH  (w1:ifc__IfcWall)-[r1:ifc__name_IfcRoot]->(l1:ifc__IfcLabel{express__hasString:["Base Wall"]}),
   (w2:ifc__IfcWall)-[r2:ifc__name_IfcRoot]->(l2:ifc__IfcLabel{express__hasString:["Test Wall 2"]}),
       (w3:ifc__IfcWall)-[r3:ifc__name_IfcRoot]->(l3:ifc__IfcLabel{express__hasString:["Test Wall 3"]}),
       (w4:ifc__IfcWall)-[r4:ifc__name_IfcRoot]->(l4:ifc__IfcLabel{express__hasString:["Test Wall 4"]}),
       (w1)<-[r5:ifc__relatingElement_IfcRelConnectsElements]- (n1:ifc__IfcRelConnectsPathElements)-[r6:ifc__relatedElement_IfcRelConnectsElements]->(w2) 
	  <-[r7:ifc__relatingElement_IfcRelConnectsElements]-(n2:ifc__IfcRelConnectsPathElements) -[r8:ifc__relatedElement_IfcRelConnectsElements]-> (w3) 
	 <-[r9:ifc__relatingElement_IfcRelConnectsElements]-(n3:ifc__IfcRelConnectsPathElements)-[r10:ifc__relatedElement_IfcRelConnectsElements]->(w4) 
             <-[r11:ifc__relatedElement_IfcRelConnectsElements]-(n4:ifc__IfcRelConnectsPathElements)-[r12:ifc__relatingElement_IfcRelConnectsElements]->(w1)
DETACH DELETE w1,w2,w3,w4,l1,l2,l3,l4,n1,n2,n3,n4,r1,r2,r3,r4,r5,r6,r7,r8,r9,r10,r11,r12
'@

$d6 = @'
MATCH (:IfcElement{Entity:"IfcWall",Name:"Base Wall"})-[r1:IfcRelVoidsElement]->(o:IfcElement{Entity:"IfcOpeningElement"})-[r2:IfcRelFillsElement]->(d:IfcElement{Entity:"IfcDoor",Name:"Base Door"})
DETACH DELETE o,d,r1,r2
'@

$e6 = @'
MATCH (l:ifc__IfcLabel{express__hasString:["Base Door"]})<-[r1:ifc__name_IfcRoot]-(n1:ifc__IfcDoor)
				<-[r2:ifc__relatedBuildingElement_IfcRelFillsElement]-(n2:ifc__IfcRelFillsElement)
				-[r3:ifc__relatingOpeningElement_IfcRelFillsElement]->(n3:ifc__IfcOpeningElement)
				<-[r4:ifc__relatedOpeningElement_IfcRelVoidsElement]-(n4:ifc__IfcRelVoidsElement)
				-[r5:ifc__relatingBuildingElement_IfcRelVoidsElement]->(:ifc__IfcWall)
DETACH DELETE l,r1,r2,r3,r4,r5,n1,n2,n3,n4
'@

$d7 = @'
MATCH (d {Entity:"IfcSite"})-[r1:IfcRelContainedInSpatialStructure]->(w1:IfcElement{Entity:"IfcWall",Name:"Base Wall"}),
      (d)-[r2:IfcRelContainedInSpatialStructure]->(w2:IfcElement{Entity:"IfcWall",Name:"Test Wall 2"}),
      (d)-[r3:IfcRelContainedInSpatialStructure]->(w3:IfcElement{Entity:"IfcWall",Name:"Test Wall 3"}),
      (d)-[r4:IfcRelContainedInSpatialStructure]->(w4:IfcElement{Entity:"IfcWall",Name:"Test Wall 4"}),
      (w1)-[r5:IfcRelConnectsPathElements]->(w2)-[r6:IfcRelConnectsPathElements]->(w3)<-[r7:IfcRelConnectsPathElements]-(w4)-[r8:IfcRelConnectsPathElements]->(w1)
DETACH DELETE w1,w2,w3,w4,r1,r2,r3,r4,r5,r6,r7,r8
'@

$ws.Range("E7").Value = $e7
$ws.Range("D6").Value = $d6
$ws.Range("E6").Value = $e6
$ws.Range("D7").Value = $d7

$ws.Range("D6:E7").WrapText = $true

$ws.Rows.Item(6).RowHeight = 102
$ws.Rows.Item(7).RowHeight = 221

$ws.Range("D18").Select()
